$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (cols C, D, F, H) ---
# Excel's ColumnWidth property (character units) differs from the raw
# OOXML <col width="..."> value by a constant padding offset - for the
# Calibri 11 default font used in this workbook that offset is 5/6 of a
# character. Subtract it here so the saved XML width lands exactly on
# the desired integer value (76 / 34 / 16 / 45).
$offset = 5.0 / 6.0
$ws.Columns.Item(3).ColumnWidth = 76 - $offset
$ws.Columns.Item(4).ColumnWidth = 34 - $offset
$ws.Columns.Item(6).ColumnWidth = 16 - $offset
$ws.Columns.Item(8).ColumnWidth = 45 - $offset

# --- Row 2: refreshed scrape entry ---
$ws.Range("A2").Value = "1332087"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1332087"
$ws.Range("C2").Value = "Administrative and commercial assistant"
$ws.Range("D2").Value = "Lisboa, Portugal"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "18 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "AAVANZ - INOVAÇÃO, UNIPESSOAL LDA"

# --- Row 3: refreshed scrape entry (also drops the old yellow "Premium"
#     highlight on column E, since this listing is no longer premium) ---
$ws.Range("A3").Value = "1332032"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1332032"
$ws.Range("C3").Value = "Sales Representative"
$ws.Range("D3").Value = "İstanbul, Türkiye"
$ws.Range("E3").ClearFormats()
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "7 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Boğaziçi İhtisas Fuarcılık Limited Şirketi"

# --- Row 4: refreshed scrape entry ---
$ws.Range("A4").Value = "1331961"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1331961"
$ws.Range("C4").Value = "Social Media Marketing"
$ws.Range("D4").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "Aman Real Estate"

# --- Row 5: refreshed scrape entry ---
$ws.Range("A5").Value = "1331509"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1331509"
$ws.Range("C5").Value = "Commercial area internship - KAM Jr."
$ws.Range("D5").Value = "Antiguo Cuscatlán, El Salvador"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "34 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "Samsung Electronics Latinoamerica"

# --- Row 6: refreshed scrape entry ---
$ws.Range("A6").Value = "1331459"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1331459"
$ws.Range("C6").Value = "Financial Accounting Trainee (EU ONLY)"
$ws.Range("D6").Value = "Amsterdam, Nederland"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "2 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Volkswagen International Finance N.V."

# --- Row 7: refreshed scrape entry ---
$ws.Range("A7").Value = "1317170"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1317170"
$ws.Range("C7").Value = "Guest Relations Officer"
$ws.Range("D7").Value = "Colombo, Sri Lanka"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "47 applicants"
$ws.Range("G7").Value = "3 - 6 Months"
$ws.Range("H7").Value = "Lanka Island Resorts Ltd"

# --- Row 8: refreshed scrape entry ---
$ws.Range("A8").Value = "1304488"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1304488"
$ws.Range("C8").Value = "Client Consultant | Tourism Sector (German speaker) ( Flexible RE dates )"
$ws.Range("D8").Value = "Athens, Greece"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "84 applicants"
$ws.Range("G8").Value = "3 - 6 Months"
$ws.Range("H8").Value = "SpeakIT"

# --- Row 9: refreshed scrape entry ---
$ws.Range("A9").Value = "1289375"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1289375"
$ws.Range("C9").Value = "Medical Advisor (German Speaker)"
$ws.Range("D9").Value = "İstanbul, Türkiye"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "45 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "International Plus"

# --- Row 10 no longer present in the latest scrape: drop it entirely,
#     which also shrinks the sheet's used range/dimension to A1:H9 ---
$ws.Rows("10:10").Delete()
